$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.719.95"
$ws.Range("E2").Value = "  -1.44%  "
$ws.Range("D3").Value = "1.804.51"
$ws.Range("E3").Value = "  -1.00%  "
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  +0.21%  "
$ws.Range("D5").Value = "232.02"
$ws.Range("E5").Value = "  -1.29%  "
$ws.Range("D6").Value = "0.5921"
$ws.Range("E6").Value = "  -2.36%  "
$ws.Range("E7").Value = "  +0.19%  "
$ws.Range("D8").Value = "0.2779"
$ws.Range("E8").Value = "  -0.53%  "
$ws.Range("D9").Value = "0.06829"
$ws.Range("E9").Value = "  -3.36%  "
$ws.Range("D10").Value = "23.34"
$ws.Range("E10").Value = "  -0.39%  "
$ws.Range("D11").Value = "0.07493"
$ws.Range("E11").Value = "  -2.01%  "
$ws.Range("D12").Value = "1.798.85"
$ws.Range("E12").Value = "  -1.17%  "
$ws.Range("D13").Value = "4.776"
$ws.Range("E13").Value = "  -0.23%  "
$ws.Range("D14").Value = "0.6233"
$ws.Range("D15").Value = "2.050.74"
$ws.Range("E15").Value = "  -0.87%  "
$ws.Range("D16").Value = "0.000009292"
$ws.Range("E16").Value = "  -6.25%  "
$ws.Range("D17").Value = "75.69"
$ws.Range("E17").Value = "  -3.53%  "
$ws.Range("D18").Value = "28.658.67"
$ws.Range("E18").Value = "  -1.59%  "
$ws.Range("D19").Value = "5.486"
$ws.Range("E19").Value = "  -6.07%  "
$ws.Range("D20").Value = "1.003"
$ws.Range("E20").Value = "  +0.20%  "
$ws.Range("D21").Value = "210.96"
$ws.Range("E21").Value = "  -6.65%  "
$ws.Range("D22").Value = "11.47"
$ws.Range("E22").Value = "  -2.16%  "
$ws.Range("D23").Value = "6.853"
$ws.Range("E23").Value = "  -1.45%  "
$ws.Range("D24").Value = "1.004"
$ws.Range("E24").Value = "  +0.46%  "
$ws.Range("D25").Value = "154.24"
$ws.Range("E25").Value = "  -0.59%  "
$ws.Range("D26").Value = "7.878"
$ws.Range("E26").Value = "  -1.65%  "
$ws.Range("D27").Value = "0.1272"
$ws.Range("E27").Value = "  -2.16%  "
$ws.Range("D28").Value = "16.43"
$ws.Range("E28").Value = "  -0.61%  "
$ws.Range("D29").Value = "1.421"
$ws.Range("E29").Value = "  -4.85%  "
$ws.Range("D30").Value = "0.06230"
$ws.Range("E30").Value = "  -0.19%  "
$ws.Range("D31").Value = "1.423"
$ws.Range("E31").Value = "  -1.54%  "
$ws.Range("D32").Value = "3.781"
$ws.Range("E32").Value = "  -1.11%  "
$ws.Range("D33").Value = "3.759"
$ws.Range("E33").Value = "  -0.77%  "
$ws.Range("D34").Value = "1.726"
$ws.Range("E34").Value = "  -0.75%  "
$ws.Range("D35").Value = "1.064"
$ws.Range("E35").Value = "  -4.99%  "
$ws.Range("D36").Value = "0.6410"
$ws.Range("E36").Value = "  +0.52%  "
$ws.Range("D37").Value = "2.492"
$ws.Range("D38").Value = "2.716"
$ws.Range("E38").Value = "  +0.03%  "
$ws.Range("D39").Value = "6.510"
$ws.Range("E39").Value = "  +0.21%  "
$ws.Range("D40").Value = "0.01716"
$ws.Range("E40").Value = "  -1.08%  "
$ws.Range("D41").Value = "1.143.58"
$ws.Range("E41").Value = "  -5.51%  "
$ws.Range("D42").Value = "0.8787"
$ws.Range("E42").Value = "  -2.67%  "
$ws.Range("D43").Value = "1.007"
$ws.Range("E43").Value = "  +0.63%  "
$ws.Range("D44").Value = "100.29"
$ws.Range("E44").Value = "  -0.22%  "
$ws.Range("D45").Value = "1.962.34"
$ws.Range("E45").Value = "  -0.82%  "
$ws.Range("D46").Value = "60.48"
$ws.Range("E46").Value = "  -3.10%  "
$ws.Range("E47").Value = "  -1.98%  "
$ws.Range("D48").Value = "1.598"
$ws.Range("E48").Value = "  +0.34%  "
$ws.Range("D49").Value = "8.414"
$ws.Range("E49").Value = "  -0.89%  "
$ws.Range("D50").Value = "0.05461"
$ws.Range("E50").Value = "  -0.78%  "
$ws.Range("D51").Value = "0.4482"
$ws.Range("E51").Value = "  -1.71%  "
